# Update "想去人数" (F column) counts across sheets to the newly scraped
# values, matching the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 14924
$ws1.Range("F3").Value  = 18667
$ws1.Range("F5").Value  = 120
$ws1.Range("F13").Value = 53
$ws1.Range("F14").Value = 120
$ws1.Range("F15").Value = 204
$ws1.Range("F16").Value = 57
$ws1.Range("F17").Value = 1430
$ws1.Range("F21").Value = 232
$ws1.Range("F22").Value = 7762
$ws1.Range("F24").Value = 23
$ws1.Range("F25").Value = 56
$ws1.Range("F26").Value = 1226
$ws1.Range("F28").Value = 5976
$ws1.Range("F30").Value = 67
$ws1.Range("F31").Value = 158
$ws1.Range("F33").Value = 265
$ws1.Range("F34").Value = 5344

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value  = 3

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 14924
$ws4.Range("F3").Value  = 18667
$ws4.Range("F5").Value  = 120
$ws4.Range("F13").Value = 53
$ws4.Range("F14").Value = 120
$ws4.Range("F15").Value = 204
$ws4.Range("F16").Value = 57
$ws4.Range("F17").Value = 1430
$ws4.Range("F22").Value = 232
$ws4.Range("F23").Value = 7762
$ws4.Range("F25").Value = 23
$ws4.Range("F26").Value = 56
$ws4.Range("F27").Value = 1226
$ws4.Range("F29").Value = 3
$ws4.Range("F31").Value = 5976
$ws4.Range("F33").Value = 67
$ws4.Range("F34").Value = 158
$ws4.Range("F35").Value = 154
$ws4.Range("F36").Value = 265
$ws4.Range("F37").Value = 5344
